$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "enum:byte"
$ws.Range("A4").Value = "New = 0"
$ws.Range("A5").Select()
